$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.848.24'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '1.890.70'
$ws.Range('E3').Value = '  -0.86%  '
$ws.Range('E4').Value = '  -0.55%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7743'
$ws.Range('E5').Value = '  -2.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.70'
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('E7').Value = '  -0.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3139'
$ws.Range('E8').Value = '  -1.99%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07365'
$ws.Range('E9').Value = '  +3.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.32'
$ws.Range('E10').Value = '  -3.87%  '
$ws.Range('E11').Value = '  +0.96%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7662'
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.477'
$ws.Range('E13').Value = '  +2.71%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.900.63'
$ws.Range('E14').Value = '  -1.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.27'
$ws.Range('E15').Value = '  -0.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.186'
$ws.Range('E16').Value = '  +3.52%  '
$ws.Range('D17').Value = '29.868.93'
$ws.Range('E17').Value = '  -0.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.95'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.05'
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007832'
$ws.Range('E20').Value = '  +1.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.175'
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('D23').Value = '2.138.36'
$ws.Range('E23').Value = '  -1.52%  '
$ws.Range('E24').Value = '  -0.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1575'
$ws.Range('E25').Value = '  -2.00%  '
$ws.Range('E26').Value = '  +0.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.38'
$ws.Range('E27').Value = '  -2.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.81'
$ws.Range('E28').Value = '  +0.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.039'
$ws.Range('E29').Value = '  -3.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.455'
$ws.Range('E30').Value = '  +5.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.548'
$ws.Range('E31').Value = '  +0.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.496'
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05614'
$ws.Range('E33').Value = '  -1.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.093'
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.248'
$ws.Range('E35').Value = '  -1.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7621'
$ws.Range('E36').Value = '  +3.19%  '
$ws.Range('E37').Value = '  +0.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.646'
$ws.Range('E38').Value = '  -3.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01926'
$ws.Range('E39').Value = '  -0.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.789'
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('D41').Value = '1.163.09'
$ws.Range('E41').Value = '  +12.91%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '74.23'
$ws.Range('E42').Value = '  +1.92%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4457'
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8516'
$ws.Range('E45').Value = '  +0.53%  '
$ws.Range('E46').Value = '  -0.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.904'
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.33'
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.919'
$ws.Range('E49').Value = '  -0.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.082'
$ws.Range('E50').Value = '  +1.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.534'
$ws.Range('E51').Value = '  +0.34%  '
